# Montenegro Prva Liga update (06-04-2024 01:36)
# - Swap the row order of several fixtures that were inserted out of
#   chronological order (same match-day pairs got their rows swapped).
# - Rotate 4 rows (130-133) that also needed reordering.
# - Append 5 new fixture rows (135-139) at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Simple pairwise row swaps (columns B:AC - column A, the running id, is
#    left untouched since it always matches the physical row number).
# ---------------------------------------------------------------------------
$pairs = @(
    @(17,18),
    @(25,26),
    @(29,30),
    @(38,39),
    @(53,54),
    @(59,60),
    @(63,64),
    @(76,77),
    @(85,86),
    @(110,111)
)

foreach ($p in $pairs) {
    $r1 = $p[0]
    $r2 = $p[1]
    $range1 = $ws.Range("B$r1" + ":AC$r1")
    $range2 = $ws.Range("B$r2" + ":AC$r2")
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value2 = $v2
    $range2.Value2 = $v1
}

# ---------------------------------------------------------------------------
# 2) Rows 130-133 rotate as a 4-cycle:
#    new130 = old133, new131 = old132, new132 = old130, new133 = old131
# ---------------------------------------------------------------------------
$v130 = $ws.Range("B130:AC130").Value2
$v131 = $ws.Range("B131:AC131").Value2
$v132 = $ws.Range("B132:AC132").Value2
$v133 = $ws.Range("B133:AC133").Value2

$ws.Range("B130:AC130").Value2 = $v133
$ws.Range("B131:AC131").Value2 = $v132
$ws.Range("B132:AC132").Value2 = $v130
$ws.Range("B133:AC133").Value2 = $v131

# ---------------------------------------------------------------------------
# 3) Append 5 new fixture rows (135-139).
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row=135; A=133; B=8043517; C="Montenegro Prva Liga"; D="Montenegro Prva Liga"; E=45385.41666666666; F="FK Jedinstvo Bijelo Polje"; G="FK Decic Tuzi"; H=0; I=2; J="A"; K=5.5; L=3.2; M=1.615; N=6.5; O=3.4; P=1.533; Q=1; R=1.85; S=1.95; T=2.25; U=2.025; V=1.775; W=-1; X=-1; Y=0.5329999999999999; Z=-1; AA=0.95; AB=-0.5; AC=0.3875 },
    @{ Row=136; A=134; B=8043518; C="Montenegro Prva Liga"; D="Montenegro Prva Liga"; E=45385.41666666666; F="FK Arsenal"; G="FK Rudar Pljevlja"; H=4; I=2; J="H"; K=1.909; L=3; M=3.9; N=1.65; O=3.3; P=5; Q=-0.75; R=1.875; S=1.925; T=2; U=1.8; V=2; W=0.6499999999999999; X=-1; Y=-1; Z=0.875; AA=-1; AB=0.8; AC=-1 },
    @{ Row=137; A=135; B=6815409; C="Montenegro Prva Liga"; D="Montenegro Prva Liga"; E=45385.5; F="Sutjeska Niksic"; G="FK Jezero"; H=2; I=1; J="H"; K=1.727; L=3; M=5; N=1.909; O=3; P=4; Q=-0.5; R=1.925; S=1.875; T=2.25; U=2; V=1.8; W=0.909; X=-1; Y=-1; Z=0.925; AA=-1; AB=1; AC=-1 },
    @{ Row=138; A=136; B=8043515; C="Montenegro Prva Liga"; D="Montenegro Prva Liga"; E=45385.54166666666; F="FK Mornar Bar"; G="OFK Petrovac"; H=1; I=0; J="H"; K=2.375; L=2.8; M=3; N=2.1; O=2.8; P=3.6; Q=-0.25; R=1.825; S=1.975; T=2; U=1.95; V=1.85; W=1.1; X=-1; Y=-1; Z=0.825; AA=-1; AB=-1; AC=0.8500000000000001 },
    @{ Row=139; A=137; B=8043516; C="Montenegro Prva Liga"; D="Montenegro Prva Liga"; E=45385.5625; F="Buducnost Podgorica"; G="OFK Mladost DG"; H=1; I=2; J="A"; K=1.3; L=4.4; M=8.5; N=1.25; O=5.25; P=9; Q=-1.75; R=1.8; S=2; T=3.25; U=1.85; V=1.95; W=-1; X=-1; Y=8; Z=-1; AA=1; AB=-0.5; AC=0.475 }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

foreach ($nr in $newRows) {
    $r = $nr.Row
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value2 = $nr[$col]
    }
}

# Copy the number-formats of an existing data row onto the new rows so the
# id column (A) and date column (E) keep the same styles (s="1" / s="2")
# instead of picking up a brand-new style index.
$ws.Range("A17:AC17").Copy()
foreach ($nr in $newRows) {
    $ws.Range("A$($nr.Row):AC$($nr.Row)").PasteSpecial(-4122)
}
